# Questionnaire Breakdown - "Worked on Q breakdown"
# Fill in the new B1 / B1a overall-consumption rows (sheet rows 26-45)
# on Sheet1, matching the Q Num rows already present in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cVals = @(
  "B1:1","B1:2","B1:3","B1:4","B1:5","B1:6","B1:7","B1:8","B1:9","B1:10","B1:11",
  "B1a: 1","B1a: 2","B1a: 3","B1a: 4","B1a: 5","B1a: 6","B1a: 7","B1a: 8","B1a: 9"
)

$dVals = @(
  "TV:Buy: Full Season: DVD (# P6M: Total)",
  "TV:Buy: Full Season: Blu-ray (# P6M: Total)",
  "TV:Buy: Full Season: Digital (# P6M: Total)",
  "TV:Buy: Episode: Digital (# P6M: Total)",
  "TV: Rent: Full Season: Disc: Walk-in (# P6M: Total)",
  "TV: Rent: Full Season: Disc: Mail (# P6M: Total)",
  "TV: Rent: Full Season: Digital: Sub Stream (# P6M: Total)",
  "TV: Rent: Episode: Digital: Sub Stream (# P6M: Total)",
  "TV: Free: Full Season: Digital: Ad Stream (# P6M: Total)",
  "TV: Free: Episode: Digital: Ad Stream (# P6M: Total)",
  "TV: Free: Episode: Digital: cVOD (# P6M: Total)",
  "TV:Buy: Full Season: DVD (# P6M: Fam)",
  "TV:Buy: Full Season: Blu-ray (# P6M: Fam)",
  "TV:Buy: Full Season: Digital (# P6M: Fam)",
  "TV:Buy: Episode: Digital (# P6M: Fam)",
  "TV: Rent: Full Season: Disc: Walk-in (# P6M: Fam)",
  "TV: Rent: Full Season: Disc: Mail (# P6M: Fam)",
  "TV: Rent: Full Season: Digital: Sub Stream (# P6M: Fam)",
  "TV: Rent: Episode: Digital: Sub Stream (# P6M: Fam)",
  "TV: Free: Full Season: Digital: Ad Stream (# P6M: Fam)"
)

for ($i = 0; $i -lt 20; $i++) {
    $r = 26 + $i
    $ws.Cells.Item($r, 2).Value = "Overall consumption"
    $ws.Cells.Item($r, 3).Value = $cVals[$i]
    $ws.Cells.Item($r, 4).Value = $dVals[$i]
    $ws.Cells.Item($r, 5).Value = "OE"
    if ($i -ge 11) {
        $ws.Cells.Item($r, 6).Value = "1+ at corresponding B1 Q"
    }
}

# Column width tweaks that accompanied the new, wider text in columns B/D,
# plus a stray widened (empty) column I.
$ws.Columns.Item(2).ColumnWidth = 19.17
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(4).ColumnWidth = 51.17
$ws.Columns.Item(9).ColumnWidth = 21.5

# View state: scrolled further down, zoomed to 85%, new selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 24
$win.ScrollColumn = 1
$win.Zoom = 85
$ws.Range("D45").Select()

Write-Output "applied Q breakdown updates"
